$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing segment names (currently in column A, rows 2-20)
# before we shift columns around.
$lastRow = 20
$names = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $names += $ws.Cells.Item($r, 1).Value()
}

# Insert a new column before column B; this shifts the old B:E -> C:F.
$ws.Columns("B").Insert()

# Give the new header cell (B1) the same formatting as the other header
# cells (bold / centered / bordered), then set its text.
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)
$ws.Cells.Item(1, 2).Value = "segments"

# Fill in the segment names (now in column B, unstyled like the numeric
# data columns) and replace column A with a numeric 0-based index.
for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $i
    $ws.Cells.Item($r, 2).Value = $names[$i]
}

$excel.CutCopyMode = 0
